$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 3800
$ws.Range("J46").Value = 10000
$ws.Range("L46").Value = 30000
$ws.Range("N46").Value = -30238
$ws.Range("H55").Value = 336.2
$ws.Range("I55").Value = 336.2
$ws.Range("K55").Value = 336.2
$ws.Range("M55").Value = -122.2
$ws.Range("H60").Value = 3800
$ws.Range("J60").Value = 10000
$ws.Range("L60").Value = 30000
$ws.Range("N60").Value = -30968
$ws.Range("H98").Value = 688.51514
$ws.Range("I98").Value = 692.6875
$ws.Range("K98").Value = 692.6875
$ws.Range("M98").Value = 805.3125
$ws.Range("H111").Value = 18333.2
$ws.Range("I111").Value = 10243.4
$ws.Range("J111").Value = 26423
$ws.Range("K111").Value = 30730.2
$ws.Range("L111").Value = 79269
$ws.Range("M111").Value = -27663.2
$ws.Range("N111").Value = -85403
$ws.Range("H115").Value = 261.42856
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("H122").Value = 688.51514
$ws.Range("I122").Value = 692.6875
$ws.Range("K122").Value = 2078.0625
$ws.Range("M122").Value = 371.9375
$ws.Range("H132").Value = 3115.4753
$ws.Range("I132").Value = 2923.698
$ws.Range("K132").Value = 8771.093999999999
$ws.Range("M132").Value = -6241.093999999999
$ws.Range("H138").Value = 2039.6774
$ws.Range("I138").Value = 1697.36
$ws.Range("K138").Value = 5092.08
$ws.Range("M138").Value = 47.92000000000007
$ws.Range("N115").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2683.75
$ws.Range("I32").Value = 2759.5454
$ws.Range("J32").Value = 1850
$ws.Range("K32").Value = 2759.5454
$ws.Range("L32").Value = 1850
$ws.Range("M32").Value = -2472.5454
$ws.Range("N32").Value = -2424
$ws.Range("H74").Value = 3108.8438
$ws.Range("I74").Value = 3049.6333
$ws.Range("J74").Value = 3997
$ws.Range("K74").Value = 3049.6333
$ws.Range("L74").Value = 3997
$ws.Range("M74").Value = -2175.6333
$ws.Range("N74").Value = -5745
$ws.Range("H77").Value = 3108.8438
$ws.Range("I77").Value = 3049.6333
$ws.Range("J77").Value = 3997
$ws.Range("K77").Value = 15248.1665
$ws.Range("L77").Value = 19985
$ws.Range("M77").Value = -10880.1665
$ws.Range("N77").Value = -28721
$ws.Range("H97").Value = 709.7143
$ws.Range("I97").Value = 710.45
$ws.Range("K97").Value = 710.45
$ws.Range("M97").Value = -214.45
$ws.Range("H102").Value = 61539.816
$ws.Range("I102").Value = 47173.11
$ws.Range("J102").Value = 126190
$ws.Range("K102").Value = 47173.11
$ws.Range("L102").Value = 126190
$ws.Range("M102").Value = -45551.11
$ws.Range("N102").Value = -129434
$ws.Range("H132").Value = 1547.8
$ws.Range("I132").Value = 1497.8
$ws.Range("K132").Value = 4493.4
$ws.Range("M132").Value = -1963.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4295.9414
$ws.Range("I94").Value = 4376.9375
$ws.Range("K94").Value = 4376.9375
$ws.Range("M94").Value = -3925.9375
$ws.Range("H99").Value = 1398.909
$ws.Range("J99").Value = 1798.3334
$ws.Range("L99").Value = 1798.3334
$ws.Range("N99").Value = -4794.3334
$ws.Range("H134").Value = 1312.6364
$ws.Range("I134").Value = 1158.9
$ws.Range("J134").Value = 2850
$ws.Range("K134").Value = 3476.7
$ws.Range("L134").Value = 8550
$ws.Range("M134").Value = -941.7000000000003
$ws.Range("N134").Value = -13620

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 67760.92999999999
$ws.Range("I86").Value = 71819.84
$ws.Range("K86").Value = 71819.84
$ws.Range("M86").Value = -70696.84
$ws.Range("H89").Value = 67760.92999999999
$ws.Range("I89").Value = 71819.84
$ws.Range("K89").Value = 359099.2
$ws.Range("M89").Value = -353483.2
$ws.Range("H99").Value = 2096.3845
$ws.Range("I99").Value = 2096.3845
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2096.3845
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -598.3845000000001
$ws.Range("H122").Value = 2652
$ws.Range("I122").Value = 2633.5715
$ws.Range("K122").Value = 7900.7145
$ws.Range("M122").Value = -5450.7145
$ws.Range("H126").Value = 2096.3845
$ws.Range("I126").Value = 2096.3845
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6289.1535
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3819.1535
$ws.Range("N99","N126").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 899.6
$ws.Range("J86").Value = 999.6667
$ws.Range("L86").Value = 2999.0001
$ws.Range("N86").Value = -5371.0001
$ws.Range("H89").Value = 899.6
$ws.Range("J89").Value = 999.6667
$ws.Range("L89").Value = 8997.0003
$ws.Range("N89").Value = -20853.0003
$ws.Range("H131").Value = 5280.727
$ws.Range("I131").Value = 831.44446
$ws.Range("J131").Value = 8361
$ws.Range("K131").Value = 2494.33338
$ws.Range("L131").Value = 25083
$ws.Range("M131").Value = 2545.66662
$ws.Range("N131").Value = -35163

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3534.3333
$ws.Range("I80").Value = 3720.375
$ws.Range("J80").Value = 3263.7273
$ws.Range("K80").Value = 3720.375
$ws.Range("L80").Value = 3263.7273
$ws.Range("M80").Value = -2722.375
$ws.Range("N80").Value = -5259.7273
$ws.Range("H83").Value = 3534.3333
$ws.Range("I83").Value = 3720.375
$ws.Range("J83").Value = 3263.7273
$ws.Range("K83").Value = 18601.875
$ws.Range("L83").Value = 16318.6365
$ws.Range("M83").Value = -13609.875
$ws.Range("N83").Value = -26302.6365
$ws.Range("H97").Value = 1646.9048
$ws.Range("I97").Value = 1693
$ws.Range("K97").Value = 1693
$ws.Range("M97").Value = -1197
$ws.Range("H102").Value = 2629.7307
$ws.Range("I102").Value = 1547.2
$ws.Range("K102").Value = 1547.2
$ws.Range("M102").Value = 74.79999999999995
$ws.Range("H132").Value = 24353.273
$ws.Range("I132").Value = 38567.418
$ws.Range("J132").Value = 7296.3
$ws.Range("K132").Value = 115702.254
$ws.Range("L132").Value = 21888.9
$ws.Range("M132").Value = -113172.254
$ws.Range("N132").Value = -26948.9

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2767.318
$ws.Range("I40").Value = 2180.125
$ws.Range("K40").Value = 2180.125
$ws.Range("M40").Value = -2044.125
$ws.Range("H46").Value = 2659.2812
$ws.Range("I46").Value = 1706.25
$ws.Range("K46").Value = 1706.25
$ws.Range("M46").Value = -1518.25
$ws.Range("H61").Value = 1268.2222
$ws.Range("I61").Value = 1172.7142
$ws.Range("J61").Value = 1602.5
$ws.Range("K61").Value = 1172.7142
$ws.Range("L61").Value = 1602.5
$ws.Range("M61").Value = -970.7141999999999
$ws.Range("N61").Value = -2006.5
$ws.Range("H113").Value = 1268.2222
$ws.Range("I113").Value = 1172.7142
$ws.Range("J113").Value = 1602.5
$ws.Range("K113").Value = 1172.7142
$ws.Range("L113").Value = 1602.5
$ws.Range("M113").Value = 997.2858000000001
$ws.Range("N113").Value = -5942.5
$ws.Range("H122").Value = 4427.923
$ws.Range("I122").Value = 3066.353
$ws.Range("J122").Value = 6999.778
$ws.Range("K122").Value = 9199.059000000001
$ws.Range("L122").Value = 20999.334
$ws.Range("M122").Value = -6749.059000000001
$ws.Range("N122").Value = -25899.334
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N124","N125","N127").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5996.6665
$ws.Range("I14").Value = 3995
$ws.Range("K14").Value = 3995
$ws.Range("M14").Value = -3827
$ws.Range("H122").Value = 1570.7941
$ws.Range("I122").Value = 1454.4
$ws.Range("J122").Value = 1894.1111
$ws.Range("K122").Value = 4363.200000000001
$ws.Range("L122").Value = 5682.3333
$ws.Range("M122").Value = -1913.200000000001
$ws.Range("N122").Value = -10582.3333
$ws.Range("H126").Value = 1956.1428
$ws.Range("I126").Value = 1956.1428
$ws.Range("K126").Value = 5868.428400000001
$ws.Range("M126").Value = -3398.428400000001
$ws.Range("H132").Value = 3611.158
$ws.Range("I132").Value = 3327.1292
$ws.Range("K132").Value = 9981.3876
$ws.Range("M132").Value = -7451.3876
